# "Repayment Schedule" sheet (Loan RBI / Variable Instalments):
# insert a new blank column before the old "Late" column (N), pushing
# Late -> O, the blank spacer -> P, and "Outstanding" -> Q.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

$ws.Columns("N:N").Insert()

# Match the width Excel copies from the left-hand neighbour (M) when a
# column is inserted, instead of the engine's computed best-fit width.
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

$ws.Range("R9").Select()
